$d = $word.ActiveDocument

$d.Content.Find.Execute("79×26=", $true, $false, $false, $false, $false, $true, 1, $false, "98×50=", 2) | Out-Null
$d.Content.Find.Execute("43×92=", $true, $false, $false, $false, $false, $true, 1, $false, "42×28=", 2) | Out-Null
$d.Content.Find.Execute("47×90=", $true, $false, $false, $false, $false, $true, 1, $false, "30×35=", 2) | Out-Null
$d.Content.Find.Execute("21×30=", $true, $false, $false, $false, $false, $true, 1, $false, "17×55=", 2) | Out-Null
$d.Content.Find.Execute("62×13=", $true, $false, $false, $false, $false, $true, 1, $false, "69×38=", 2) | Out-Null
$d.Content.Find.Execute("61×54=", $true, $false, $false, $false, $false, $true, 1, $false, "24×71=", 2) | Out-Null
$d.Content.Find.Execute("13×36=", $true, $false, $false, $false, $false, $true, 1, $false, "11×92=", 2) | Out-Null
$d.Content.Find.Execute("54×76=", $true, $false, $false, $false, $false, $true, 1, $false, "60×31=", 2) | Out-Null
$d.Content.Find.Execute("68×41=", $true, $false, $false, $false, $false, $true, 1, $false, "34×72=", 2) | Out-Null
$d.Content.Find.Execute("34×62=", $true, $false, $false, $false, $false, $true, 1, $false, "45×46=", 2) | Out-Null
$d.Content.Find.Execute("15×66=", $true, $false, $false, $false, $false, $true, 1, $false, "45×73=", 2) | Out-Null
$d.Content.Find.Execute("41×53=", $true, $false, $false, $false, $false, $true, 1, $false, "38×74=", 2) | Out-Null
$d.Content.Find.Execute("85×42=", $true, $false, $false, $false, $false, $true, 1, $false, "29×85=", 2) | Out-Null
$d.Content.Find.Execute("71×74=", $true, $false, $false, $false, $false, $true, 1, $false, "17×44=", 2) | Out-Null
$d.Content.Find.Execute("38×68=", $true, $false, $false, $false, $false, $true, 1, $false, "47×24=", 2) | Out-Null
$d.Content.Find.Execute("95×52=", $true, $false, $false, $false, $false, $true, 1, $false, "55×48=", 2) | Out-Null
$d.Content.Find.Execute("53×90=", $true, $false, $false, $false, $false, $true, 1, $false, "79×75=", 2) | Out-Null
$d.Content.Find.Execute("91×74=", $true, $false, $false, $false, $false, $true, 1, $false, "12×79=", 2) | Out-Null
$d.Content.Find.Execute("40×88=", $true, $false, $false, $false, $false, $true, 1, $false, "62×12=", 2) | Out-Null
$d.Content.Find.Execute("38×49=", $true, $false, $false, $false, $false, $true, 1, $false, "53×99=", 2) | Out-Null
$d.Content.Find.Execute("45×78=", $true, $false, $false, $false, $false, $true, 1, $false, "22×44=", 2) | Out-Null
$d.Content.Find.Execute("75×80=", $true, $false, $false, $false, $false, $true, 1, $false, "76×47=", 2) | Out-Null
$d.Content.Find.Execute("48×28=", $true, $false, $false, $false, $false, $true, 1, $false, "86×59=", 2) | Out-Null
$d.Content.Find.Execute("45×18=", $true, $false, $false, $false, $false, $true, 1, $false, "18×85=", 2) | Out-Null
$d.Content.Find.Execute("72×81=", $true, $false, $false, $false, $false, $true, 1, $false, "60×83=", 2) | Out-Null
